$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column C entirely, shifting remaining columns left
$ws.Range("C1").EntireColumn.Delete()

# Drop row 3 entirely, shifting remaining rows up
$ws.Range("A3").EntireRow.Delete()

# B1 changes from 0 to 1
$ws.Range("B1").Value = 1
